$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1828.3549
$ws.Range("J17").Value = 1828.3549
$ws.Range("L17").Value = 5485.0647
$ws.Range("N17").Value = -5821.0647
$ws.Range("H19").Value = 275.2
$ws.Range("I19").Value = 256
$ws.Range("J19").Value = 297.14285
$ws.Range("K19").Value = 256
$ws.Range("L19").Value = 297.14285
$ws.Range("M19").Value = -81
$ws.Range("N19").Value = -647.14285
$ws.Range("H33").Value = 144.61539
$ws.Range("I33").Value = 108
$ws.Range("K33").Value = 108
$ws.Range("M33").Value = 121
$ws.Range("H64").Value = 26318092
$ws.Range("I64").Value = 33335216
$ws.Range("K64").Value = 33335216
$ws.Range("M64").Value = -33334968
$ws.Range("H67").Value = 26318092
$ws.Range("I67").Value = 33335216
$ws.Range("K67").Value = 33335216
$ws.Range("M67").Value = -33334358
$ws.Range("H138").Value = 1709.6383
$ws.Range("I138").Value = 922.7353000000001
$ws.Range("J138").Value = 3767.6924
$ws.Range("K138").Value = 2768.2059
$ws.Range("L138").Value = 11303.0772
$ws.Range("M138").Value = 2371.7941
$ws.Range("N138").Value = -21583.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1208.7941
$ws.Range("I45").Value = 1234.25
$ws.Range("J45").Value = 1172.4286
$ws.Range("K45").Value = 1234.25
$ws.Range("L45").Value = 1172.4286
$ws.Range("M45").Value = -857.25
$ws.Range("N45").Value = -1926.4286
$ws.Range("H74").Value = 20819.2
$ws.Range("I74").Value = 700.6667
$ws.Range("J74").Value = 50997
$ws.Range("K74").Value = 700.6667
$ws.Range("L74").Value = 50997
$ws.Range("M74").Value = 173.3333
$ws.Range("N74").Value = -52745
$ws.Range("H77").Value = 20819.2
$ws.Range("I77").Value = 700.6667
$ws.Range("J77").Value = 50997
$ws.Range("K77").Value = 3503.3335
$ws.Range("L77").Value = 254985
$ws.Range("M77").Value = 864.6665000000003
$ws.Range("N77").Value = -263721
$ws.Range("H102").Value = 1036.7
$ws.Range("I102").Value = 985.2222
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 985.2222
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 636.7778
$ws.Range("N102").Value = -4744
$ws.Range("H132").Value = 5113.635
$ws.Range("I132").Value = 3129.5
$ws.Range("K132").Value = 9388.5
$ws.Range("M132").Value = -6858.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1875
$ws.Range("I105").Value = 1329
$ws.Range("J105").Value = 2421
$ws.Range("K105").Value = 1329
$ws.Range("L105").Value = 2421
$ws.Range("M105").Value = 418
$ws.Range("N105").Value = -5915
$ws.Range("H107").Value = 1977.5186
$ws.Range("I107").Value = 1441
$ws.Range("J107").Value = 2889.6
$ws.Range("K107").Value = 1441
$ws.Range("L107").Value = 2889.6
$ws.Range("M107").Value = 479
$ws.Range("N107").Value = -6729.6
$ws.Range("H133").Value = 59750
$ws.Range("J133").Value = 59750
$ws.Range("L133").Value = 59750
$ws.Range("N133").Value = -69870
$ws.Range("H134").Value = 1329
$ws.Range("I134").Value = 1095.409
$ws.Range("J134").Value = 2063.1428
$ws.Range("K134").Value = 3286.227
$ws.Range("L134").Value = 6189.428400000001
$ws.Range("M134").Value = -751.2270000000003
$ws.Range("N134").Value = -11259.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1101.6364
$ws.Range("I22").Value = 550.8570999999999
$ws.Range("J22").Value = 2065.5
$ws.Range("K22").Value = 550.8570999999999
$ws.Range("L22").Value = 2065.5
$ws.Range("M22").Value = -200.8570999999999
$ws.Range("N22").Value = -2765.5
$ws.Range("H31").Value = 15736.939
$ws.Range("I31").Value = 1056.4706
$ws.Range("J31").Value = 19518.273
$ws.Range("K31").Value = 1056.4706
$ws.Range("L31").Value = 19518.273
$ws.Range("M31").Value = -761.4706000000001
$ws.Range("N31").Value = -20108.273
$ws.Range("H34").Value = 15736.939
$ws.Range("I34").Value = 1056.4706
$ws.Range("J34").Value = 19518.273
$ws.Range("K34").Value = 1056.4706
$ws.Range("L34").Value = 19518.273
$ws.Range("M34").Value = -854.4706000000001
$ws.Range("N34").Value = -19922.273
$ws.Range("H35").Value = 1229.6666
$ws.Range("I35").Value = 767.375
$ws.Range("J35").Value = 4928
$ws.Range("K35").Value = 767.375
$ws.Range("L35").Value = 4928
$ws.Range("M35").Value = -473.375
$ws.Range("N35").Value = -5516
$ws.Range("H58").Value = 1869.762
$ws.Range("I58").Value = 603.1818
$ws.Range("K58").Value = 603.1818
$ws.Range("M58").Value = -400.1818
$ws.Range("H68").Value = 20400
$ws.Range("I68").Value = 13500
$ws.Range("J68").Value = 22125
$ws.Range("K68").Value = 13500
$ws.Range("L68").Value = 22125
$ws.Range("M68").Value = -12751
$ws.Range("N68").Value = -23623
$ws.Range("H71").Value = 20400
$ws.Range("I71").Value = 13500
$ws.Range("J71").Value = 22125
$ws.Range("K71").Value = 40500
$ws.Range("L71").Value = 66375
$ws.Range("M71").Value = -36756
$ws.Range("N71").Value = -73863
$ws.Range("H99").Value = 3217.6
$ws.Range("I99").Value = 1706
$ws.Range("J99").Value = 5485
$ws.Range("K99").Value = 1706
$ws.Range("L99").Value = 5485
$ws.Range("M99").Value = -208
$ws.Range("N99").Value = -8481
$ws.Range("H122").Value = 831.63635
$ws.Range("I122").Value = 469.6
$ws.Range("J122").Value = 1133.3334
$ws.Range("K122").Value = 1408.8
$ws.Range("L122").Value = 3400.0002
$ws.Range("M122").Value = 1041.2
$ws.Range("N122").Value = -8300.0002
$ws.Range("H126").Value = 3217.6
$ws.Range("I126").Value = 1706
$ws.Range("J126").Value = 5485
$ws.Range("K126").Value = 5118
$ws.Range("L126").Value = 16455
$ws.Range("M126").Value = -2648
$ws.Range("N126").Value = -21395
$ws.Range("H136").Value = 1869.762
$ws.Range("I136").Value = 603.1818
$ws.Range("K136").Value = 1809.5454
$ws.Range("M136").Value = 740.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1060.375
$ws.Range("I122").Value = 990.4286
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 2971.2858
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -521.2857999999997
$ws.Range("N122").Value = -9550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 32000
$ws.Range("J36").Value = 32000
$ws.Range("L36").Value = 32000
$ws.Range("N36").Value = -33124
$ws.Range("H68").Value = 1733.7894
$ws.Range("I68").Value = 1708.875
$ws.Range("J68").Value = 1866.6666
$ws.Range("K68").Value = 1708.875
$ws.Range("L68").Value = 1866.6666
$ws.Range("M68").Value = -959.875
$ws.Range("N68").Value = -3364.6666
$ws.Range("H71").Value = 1733.7894
$ws.Range("I71").Value = 1708.875
$ws.Range("J71").Value = 1866.6666
$ws.Range("K71").Value = 8544.375
$ws.Range("L71").Value = 9333.333000000001
$ws.Range("M71").Value = -4800.375
$ws.Range("N71").Value = -16821.333
$ws.Range("H82").Value = 1133.9231
$ws.Range("I82").Value = 1451
$ws.Range("J82").Value = 1076.2727
$ws.Range("K82").Value = 1451
$ws.Range("L82").Value = 1076.2727
$ws.Range("M82").Value = -1090
$ws.Range("N82").Value = -1798.2727
$ws.Range("H85").Value = 1133.9231
$ws.Range("I85").Value = 1451
$ws.Range("J85").Value = 1076.2727
$ws.Range("K85").Value = 1451
$ws.Range("L85").Value = 1076.2727
$ws.Range("M85").Value = -203
$ws.Range("N85").Value = -3572.2727
$ws.Range("H132").Value = 41492.668
$ws.Range("I132").Value = 41033.332
$ws.Range("J132").Value = 41645.777
$ws.Range("K132").Value = 123099.996
$ws.Range("L132").Value = 124937.331
$ws.Range("M132").Value = -120569.996
$ws.Range("N132").Value = -129997.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 931.3333
$ws.Range("I126").Value = 824.5454999999999
$ws.Range("K126").Value = 2473.6365
$ws.Range("M126").Value = -3.636499999999614
